$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename Sheet2 -> DATA
# ------------------------------------------------------------------
$wsRun = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("Sheet2")
$wsData.Name = "DATA"

# ------------------------------------------------------------------
# 2. RUNMANAGER (sheet1): add row 4
# ------------------------------------------------------------------
$wsRun.Range("A4").Value = "amazonHamburgerMenuTest"
$wsRun.Range("B4").Value = "To check whether Amazon website is working or not"
$wsRun.Range("C4").Value = "yes"
$wsRun.Range("D4").Value = "'1"
$wsRun.Range("E4").Value = "'1"

# Column A width (best-fit driven by the new, longer values)
$wsRun.Columns.Item(1).ColumnWidth = 26.25

# ------------------------------------------------------------------
# 3. DATA (sheet2): populate the data-provider table
# ------------------------------------------------------------------
# Header row
$wsData.Range("A1").Value = "testcasename"
$wsData.Range("B1").Value = "execute"
$wsData.Range("C1").Value = "browser"
$wsData.Range("D1").Value = "username"
$wsData.Range("E1").Value = "password"
$wsData.Range("F1").Value = "name"
$wsData.Range("G1").Value = "menuItem"

# Row 2
$wsData.Range("A2").Value = "loginLogoutTest"
$wsData.Range("B2").Value = "yes"
$wsData.Range("C2").Value = "chrome"
$wsData.Range("D2").Value = "Admin"
$wsData.Range("E2").Value = "admin123"
$wsData.Range("F2").Value = "Niyaz"
$wsData.Range("G2").Value = "'"

# Row 3
$wsData.Range("A3").Value = "loginLogoutTest"
$wsData.Range("B3").Value = "yes"
$wsData.Range("C3").Value = "edge"
$wsData.Range("D3").Value = "Admin"
$wsData.Range("E3").Value = "admin123"
$wsData.Range("F3").Value = "SeleniumAutomation"
$wsData.Range("G3").Value = "'"

# Row 4
$wsData.Range("A4").Value = "newTest"
$wsData.Range("B4").Value = "yes"
$wsData.Range("C4").Value = "chrome"
$wsData.Range("D4").Value = "Admin"
$wsData.Range("E4").Value = "admin123"
$wsData.Range("F4").Value = "'"
$wsData.Range("G4").Value = "'"

# Row 5
$wsData.Range("A5").Value = "newTest"
$wsData.Range("B5").Value = "yes"
$wsData.Range("C5").Value = "edge"
$wsData.Range("D5").Value = "Admin"
$wsData.Range("E5").Value = "admin123"
$wsData.Range("F5").Value = "'"
$wsData.Range("G5").Value = "'"

# Row 6
$wsData.Range("A6").Value = "loginLogoutTest"
$wsData.Range("B6").Value = "yes"
$wsData.Range("C6").Value = "edge"
$wsData.Range("D6").Value = "admin123"
$wsData.Range("E6").Value = "admin123"
$wsData.Range("F6").Value = "Subscribe"
$wsData.Range("G6").Value = "'"

# Row 7
$wsData.Range("A7").Value = "amazonHamburgerMenuTest"
$wsData.Range("B7").Value = "yes"
$wsData.Range("C7").Value = "chrome"
$wsData.Range("D7").Value = "'"
$wsData.Range("E7").Value = "'"
$wsData.Range("F7").Value = "'"
$wsData.Range("G7").Value = "Mobiles, Computers"

# Column widths on DATA sheet
$wsData.Columns.Item(1).ColumnWidth = 26.25
$wsData.Columns.Item(4).ColumnWidth = 20
$wsData.Columns.Item(5).ColumnWidth = 14.17
$wsData.Columns.Item(6).ColumnWidth = 19.8
$wsData.Columns.Item(7).ColumnWidth = 18.3

# ------------------------------------------------------------------
# 4. Selections: set DATA's stored selection first, then come back
#    to RUNMANAGER last so it remains the active tab/sheet.
# ------------------------------------------------------------------
$wsData.Range("A7").Select()
$wsRun.Range("A4").Select()

Write-Output "done"
